$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("APPLE", "BAT", "CAT", "DOG", "ELEPHANT", "FAN", "GOAT", "HILL", "ICECREAM", "JOKER")

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

$ws.Range("A11").Select()
